# delivery_orders.xlsx: change schema and add dim product to ssis
#
# The sheet is renamed from "Sheet1" to "Sheet2" and its internal sheetId
# advances from 1 to 2 (as happens in real Excel when a sheet is replaced
# by a freshly-added one). We get both effects - new sheetId *and* kept
# data/formatting - by duplicating the existing sheet (a true copy carries
# over cell data, column formatting, phoneticPr, etc.), deleting the
# original, and renaming the copy.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Sheet1")

$ws.Copy($null, $ws) | Out-Null
$wb.Worksheets("Sheet1").Delete() | Out-Null

$new = $wb.Worksheets("Sheet1 (2)")
$new.Name = "Sheet2"
$new.Select() | Out-Null

# Selected cell moves from C3 to O17.
$new.Range("O17").Select() | Out-Null

# Column A/B widths tighten slightly (A: 8.33 -> 8, B: 12.5 -> 12.33).
$new.Columns(1).ColumnWidth = 7.17
$new.Columns(2).ColumnWidth = 11.5

# Sheet view zoom resets from 163% back to the default 100%.
$excel.ActiveWindow.Zoom = 100
